$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert three new rows (55, 56, 57), each duplicated from row 54 so
# that number formats / styles (s="2", s="5", s="3", s="6" ...) match
# exactly what Excel would produce when a user duplicates a row and edits
# it. Rows are inserted one at a time - inserting several rows at once from
# a single-row copy source does not replicate every column correctly.
$ws.Rows("54:54").Copy()
$ws.Rows("55:55").Insert(-4121)
$ws.Rows("54:54").Copy()
$ws.Rows("56:56").Insert(-4121)
$ws.Rows("54:54").Copy()
$ws.Rows("57:57").Insert(-4121)
$excel.CutCopyMode = 0

# Column A already carries the right text ("6.4.2020", the same shared
# string used by row 54) after the row duplication above, so it is left
# untouched - re-typing it would risk Excel auto-converting the text into
# a real date serial number.

# --- From / To times
$ws.Range("B55").Value = 0.58333333333333337
$ws.Range("C55").Value = 0.58680555555555558

$ws.Range("B56").Value = 0.58680555555555558
$ws.Range("C56").Value = 0.59027777777777779

$ws.Range("B57").Value = 0.59027777777777779
$ws.Range("C57").Value = 0.63194444444444442

# --- Duration formula (Time = To - From), same pattern as the rest of the
# column; assigning the whole range at once lets Excel share the formula
# across the three new rows.
$ws.Range("D55:D57").Formula = "=C55-B55"

# --- Unit (E) / Task (F) / Notes (G)
$ws.Range("E55").Value = "ROM 1"
$ws.Range("F55").Value = "Improve Code"
$ws.Range("G55").Value = "Add packages"

$ws.Range("E56").Value = "ROM 2"
$ws.Range("F56").Value = "Improve Code"
$ws.Range("G56").Value = "Add packages"

$ws.Range("E57").Value = "VGA Top"
$ws.Range("F57").Value = "Improve Code"
$ws.Range("G57").Value = "Add packages"

# --- Match the saved selection / scroll position from the authored file.
$ws.Range("G57").Select()
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1
